$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2) ---
$ws.Range("D2").Value = "NomModule"
$ws.Range("E2").Value = "Ensiegnant_Email"
$ws.Range("F2").Value = "Classes"

# --- Email column (E3:E7) : new distinct teacher addresses ---
$ws.Range("E3").Value = "profmail1@gmail.com"
$ws.Range("E4").Value = "profmail2@gmail.com"
$ws.Range("E5").Value = "profmail3@gmail.com"
$ws.Range("E6").Value = "profmail4@gmail.com"
$ws.Range("E7").Value = "profmail5@gmail.com"

# --- Classes column (F column) : fix slash direction / casing ---
$ws.Range("F4").Value = "4.Ginfo/4GTR"
$ws.Range("F5").Value = "4.Ginfo/4GTR"

# --- Remove the hyperlink that used to sit on E3 (maryem@gmail.com) ---
$e3Links = @()
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$E$3') {
        $e3Links += $h
    }
}
foreach ($h in $e3Links) {
    $h.Delete()
}

# --- Point the remaining hyperlinks (E4:E7) at the new addresses ---
foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$E$4') { $h.Address = "mailto:profmail2@gmail.com" }
    elseif ($addr -eq '$E$5') { $h.Address = "mailto:profmail3@gmail.com" }
    elseif ($addr -eq '$E$6') { $h.Address = "mailto:profmail4@gmail.com" }
    elseif ($addr -eq '$E$7') { $h.Address = "mailto:profmail5@gmail.com" }
}

# --- Column widths ---
$ws.Columns.Item(4).ColumnWidth = 13.333333333333332
$ws.Columns.Item(5).ColumnWidth = 22.0
$ws.Columns.Item(6).ColumnWidth = 18.833333333333336

# --- Selection moves to F14 ---
$ws.Range("F14").Select()
